$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt1"
$ws.Cells.Item(2, 3).Value = "Fzd2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.01948966666666667
$ws.Cells.Item(2, 8).Value = 0.058469
$ws.Cells.Item(2, 9).Value = 0.0709606244933031
$ws.Cells.Item(2, 10).Value = 0.0709606244933031
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.139245
$ws.Cells.Item(2, 14).Value = 0.417735
$ws.Cells.Item(2, 15).Value = 0.01212793695325064
$ws.Cells.Item(2, 16).Value = 0.01283499108585158
$ws.Cells.Item(2, 17).Value = 0.002713838634999999
$ws.Cells.Item(2, 18).Value = 0.024424547715
$ws.Cells.Item(2, 19).Value = 0.0008606059800180735
$ws.Cells.Item(2, 20).Value = 0.0009107789828180069

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt1"
$ws.Cells.Item(3, 3).Value = "Fzd2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01948966666666667
$ws.Cells.Item(3, 8).Value = 0.058469
$ws.Cells.Item(3, 9).Value = 0.0709606244933031
$ws.Cells.Item(3, 10).Value = 0.0709606244933031
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.380691666666667
$ws.Cells.Item(3, 14).Value = 28.142075
$ws.Cells.Item(3, 15).Value = 0.8170378621222814
$ws.Cells.Item(3, 16).Value = 0.8646708601442703
$ws.Cells.Item(3, 17).Value = 0.1828265536861111
$ws.Cells.Item(3, 18).Value = 1.645438983175
$ws.Cells.Item(3, 19).Value = 0.05797751693087037
$ws.Cells.Item(3, 20).Value = 0.06135758421699897

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt1"
$ws.Cells.Item(4, 3).Value = "Fzd2"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.01948966666666667
$ws.Cells.Item(4, 8).Value = 0.058469
$ws.Cells.Item(4, 9).Value = 0.0709606244933031
$ws.Cells.Item(4, 10).Value = 0.0709606244933031
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.006356333333333333
$ws.Cells.Item(4, 14).Value = 0.019069
$ws.Cells.Item(4, 15).Value = 0.0005536228225107701
$ws.Cells.Item(4, 16).Value = 0.0005858988234553099
$ws.Cells.Item(4, 17).Value = 0.0001238828178888889
$ws.Cells.Item(4, 18).Value = 0.001114945361
$ws.Cells.Item(4, 19).Value = 0.00003928542121910936
$ws.Cells.Item(4, 20).Value = 0.00004157574640228034

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt1"
$ws.Cells.Item(5, 3).Value = "Fzd2"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.01948966666666667
$ws.Cells.Item(5, 8).Value = 0.058469
$ws.Cells.Item(5, 9).Value = 0.0709606244933031
$ws.Cells.Item(5, 10).Value = 0.0709606244933031
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.05759666666666666
$ws.Cells.Item(5, 14).Value = 0.17279
$ws.Cells.Item(5, 15).Value = 0.005016544522609259
$ws.Cells.Item(5, 16).Value = 0.005309007168957103
$ws.Cells.Item(5, 17).Value = 0.001122539834444444
$ws.Cells.Item(5, 18).Value = 0.01010285851
$ws.Cells.Item(5, 19).Value = 0.0003559771321228121
$ws.Cells.Item(5, 20).Value = 0.0003767304641486191

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Wnt1"
$ws.Cells.Item(6, 3).Value = "Fzd2"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.01948966666666667
$ws.Cells.Item(6, 8).Value = 0.058469
$ws.Cells.Item(6, 9).Value = 0.0709606244933031
$ws.Cells.Item(6, 10).Value = 0.0709606244933031
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.897453
$ws.Cells.Item(6, 14).Value = 3.794906
$ws.Cells.Item(6, 15).Value = 0.1652640335793479
$ws.Cells.Item(6, 16).Value = 0.1165992427774658
$ws.Cells.Item(6, 17).Value = 0.03698072648566667
$ws.Cells.Item(6, 18).Value = 0.221884358914
$ws.Cells.Item(6, 19).Value = 0.01172723902907274
$ws.Cells.Item(6, 20).Value = 0.008273955082935238

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt1"
$ws.Cells.Item(7, 3).Value = "Fzd2"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.255165
$ws.Cells.Item(7, 8).Value = 0.7654949999999999
$ws.Cells.Item(7, 9).Value = 0.9290393755066968
$ws.Cells.Item(7, 10).Value = 0.9290393755066969
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.139245
$ws.Cells.Item(7, 14).Value = 0.417735
$ws.Cells.Item(7, 15).Value = 0.01212793695325064
$ws.Cells.Item(7, 16).Value = 0.01283499108585158
$ws.Cells.Item(7, 17).Value = 0.03553045042499999
$ws.Cells.Item(7, 18).Value = 0.319774053825
$ws.Cells.Item(7, 19).Value = 0.01126733097323257
$ws.Cells.Item(7, 20).Value = 0.01192421210303358

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt1"
$ws.Cells.Item(8, 3).Value = "Fzd2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.255165
$ws.Cells.Item(8, 8).Value = 0.7654949999999999
$ws.Cells.Item(8, 9).Value = 0.9290393755066968
$ws.Cells.Item(8, 10).Value = 0.9290393755066969
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 9.380691666666667
$ws.Cells.Item(8, 14).Value = 28.142075
$ws.Cells.Item(8, 15).Value = 0.8170378621222814
$ws.Cells.Item(8, 16).Value = 0.8646708601442703
$ws.Cells.Item(8, 17).Value = 2.393624189125
$ws.Cells.Item(8, 18).Value = 21.542617702125
$ws.Cells.Item(8, 19).Value = 0.7590603451914109
$ws.Cells.Item(8, 20).Value = 0.8033132759272713

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt1"
$ws.Cells.Item(9, 3).Value = "Fzd2"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.255165
$ws.Cells.Item(9, 8).Value = 0.7654949999999999
$ws.Cells.Item(9, 9).Value = 0.9290393755066968
$ws.Cells.Item(9, 10).Value = 0.9290393755066969
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.006356333333333333
$ws.Cells.Item(9, 14).Value = 0.019069
$ws.Cells.Item(9, 15).Value = 0.0005536228225107701
$ws.Cells.Item(9, 16).Value = 0.0005858988234553099
$ws.Cells.Item(9, 17).Value = 0.001621913795
$ws.Cells.Item(9, 18).Value = 0.014597224155
$ws.Cells.Item(9, 19).Value = 0.0005143374012916607
$ws.Cells.Item(9, 20).Value = 0.0005443230770530295

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Wnt1"
$ws.Cells.Item(10, 3).Value = "Fzd2"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.255165
$ws.Cells.Item(10, 8).Value = 0.7654949999999999
$ws.Cells.Item(10, 9).Value = 0.9290393755066968
$ws.Cells.Item(10, 10).Value = 0.9290393755066969
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.05759666666666666
$ws.Cells.Item(10, 14).Value = 0.17279
$ws.Cells.Item(10, 15).Value = 0.005016544522609259
$ws.Cells.Item(10, 16).Value = 0.005309007168957103
$ws.Cells.Item(10, 17).Value = 0.01469665345
$ws.Cells.Item(10, 18).Value = 0.13226988105
$ws.Cells.Item(10, 19).Value = 0.004660567390486447
$ws.Cells.Item(10, 20).Value = 0.004932276704808483

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Wnt1"
$ws.Cells.Item(11, 3).Value = "Fzd2"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.255165
$ws.Cells.Item(11, 8).Value = 0.7654949999999999
$ws.Cells.Item(11, 9).Value = 0.9290393755066968
$ws.Cells.Item(11, 10).Value = 0.9290393755066969
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.897453
$ws.Cells.Item(11, 14).Value = 3.794906
$ws.Cells.Item(11, 15).Value = 0.1652640335793479
$ws.Cells.Item(11, 16).Value = 0.1165992427774658
$ws.Cells.Item(11, 17).Value = 0.484163594745
$ws.Cells.Item(11, 18).Value = 2.90498156847
$ws.Cells.Item(11, 19).Value = 0.1535367945502752
$ws.Cells.Item(11, 20).Value = 0.1083252876945306
